# Auto-generated Excel COM-interop script applying the Diabolos_Profits edits
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 244531.56
$ws.Range("J17").Value = 244531.56
$ws.Range("L17").Value = 733594.6799999999
$ws.Range("N17").Value = -733930.6799999999
# Row 33
$ws.Range("H33").Value = 643.2632
$ws.Range("I33").Value = 349.6
$ws.Range("J33").Value = 1744.5
$ws.Range("K33").Value = 349.6
$ws.Range("L33").Value = 1744.5
$ws.Range("M33").Value = -120.6
$ws.Range("N33").Value = -2202.5
# Row 74
$ws.Range("H74").Value = 4499.8335
$ws.Range("I74").Value = 3999
$ws.Range("K74").Value = 3999
$ws.Range("M74").Value = -3063
# Row 77
$ws.Range("H77").Value = 4499.8335
$ws.Range("I77").Value = 3999
$ws.Range("K77").Value = 19995
$ws.Range("M77").Value = -15315
# Row 80
$ws.Range("H80").Value = 1204611.6
$ws.Range("I80").Value = 1633547.9
$ws.Range("K80").Value = 4900643.699999999
$ws.Range("M80").Value = -4899645.699999999
# Row 83
$ws.Range("H83").Value = 1204611.6
$ws.Range("I83").Value = 1633547.9
$ws.Range("K83").Value = 14701931.1
$ws.Range("M83").Value = -14696939.1
# Row 113
$ws.Range("H113").Value = 125003660
$ws.Range("J113").Value = 4100
$ws.Range("L113").Value = 4100
$ws.Range("N113").Value = -10608
$ws = $wb.Worksheets.Item("ARM")
# Row 33
$ws.Range("H33").Value = 7837.6665
$ws.Range("I33").Value = 7837.6665
$ws.Range("K33").Value = 7837.6665
$ws.Range("M33").Value = -7508.6665
# Row 122
$ws.Range("H122").Value = 47621764
$ws.Range("I122").Value = 55557892
$ws.Range("K122").Value = 166673676
$ws.Range("M122").Value = -166671226
# Row 132
$ws.Range("H132").Value = 1705.15
$ws.Range("J132").Value = 798
$ws.Range("L132").Value = 2394
$ws.Range("N132").Value = -7454
$ws = $wb.Worksheets.Item("BSM")
# Row 64
$ws.Range("H64").Value = 27249.75
$ws.Range("I64").Value = 1100
$ws.Range("K64").Value = 1100
$ws.Range("M64").Value = -875
# Row 67
$ws.Range("H67").Value = 27249.75
$ws.Range("I67").Value = 1100
$ws.Range("K67").Value = 1100
$ws.Range("M67").Value = -320
# Row 99
$ws.Range("H99").Value = 1384.2106
$ws.Range("I99").Value = 1035.2858
$ws.Range("J99").Value = 2361.2
$ws.Range("K99").Value = 1035.2858
$ws.Range("L99").Value = 2361.2
$ws.Range("M99").Value = 462.7141999999999
$ws.Range("N99").Value = -5357.2
# Row 105
$ws.Range("H105").Value = 2073.077
$ws.Range("I105").Value = 1995.4546
$ws.Range("J105").Value = 2500
$ws.Range("K105").Value = 1995.4546
$ws.Range("L105").Value = 2500
$ws.Range("M105").Value = -248.4546
$ws.Range("N105").Value = -5994
$ws = $wb.Worksheets.Item("CRP")
# Row 122
$ws.Range("H122").Value = 2341.0715
$ws.Range("I122").Value = 2059.6155
$ws.Range("K122").Value = 6178.8465
$ws.Range("M122").Value = -3728.8465
$ws = $wb.Worksheets.Item("CUL")
# Row 22
$ws.Range("H22").Value = 575
$ws.Range("J22").Value = 400
$ws.Range("L22").Value = 1200
$ws.Range("N22").Value = -1538
# Row 27
$ws.Range("H27").Value = 575
$ws.Range("J27").Value = 400
$ws.Range("L27").Value = 1200
$ws.Range("N27").Value = -1404
# Row 68
$ws.Range("H68").Value = 1220.625
$ws.Range("I68").Value = 661.8333
$ws.Range("K68").Value = 1985.4999
$ws.Range("M68").Value = -1174.4999
# Row 71
$ws.Range("H71").Value = 1220.625
$ws.Range("I71").Value = 661.8333
$ws.Range("K71").Value = 5956.4997
$ws.Range("M71").Value = -1900.4997
# Row 113
$ws.Range("H113").Value = 1050.2693
$ws.Range("J113").Value = 1159
$ws.Range("L113").Value = 3477
$ws.Range("N113").Value = -7817
# Row 134
$ws.Range("H134").Value = 1356.7142
$ws.Range("I134").Value = 924.55
$ws.Range("K134").Value = 2773.65
$ws.Range("M134").Value = 2296.35
$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 991.8570999999999
$ws.Range("I97").Value = 589.1818
$ws.Range("K97").Value = 589.1818
$ws.Range("M97").Value = -93.18179999999995
# Row 122
$ws.Range("H122").Value = 2918.3572
$ws.Range("I122").Value = 2471.4285
$ws.Range("J122").Value = 3365.2856
$ws.Range("K122").Value = 7414.2855
$ws.Range("L122").Value = 10095.8568
$ws.Range("M122").Value = -4964.2855
$ws.Range("N122").Value = -14995.8568
# Row 132
$ws.Range("H132").Value = 4084.2812
$ws.Range("I132").Value = 3189.36
$ws.Range("J132").Value = 7280.4287
$ws.Range("K132").Value = 9568.08
$ws.Range("L132").Value = 21841.2861
$ws.Range("M132").Value = -7038.08
$ws.Range("N132").Value = -26901.2861
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 35716680
$ws.Range("I7").Value = 45456750
$ws.Range("J7").Value = 3093.3333
$ws.Range("K7").Value = 45456750
$ws.Range("L7").Value = 3093.3333
$ws.Range("M7").Value = -45456638
$ws.Range("N7").Value = -3317.3333
# Row 54
$ws.Range("H54").Value = 44199
$ws.Range("J54").Value = 43998.332
$ws.Range("L54").Value = 43998.332
$ws.Range("N54").Value = -45286.332
# Row 82
$ws.Range("H82").Value = 3374.5
$ws.Range("I82").Value = 3374.5
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 3374.5
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -3013.5
$ws.Range("N82").ClearContents()
# Row 85
$ws.Range("H85").Value = 3374.5
$ws.Range("I85").Value = 3374.5
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 3374.5
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -2126.5
$ws.Range("N85").ClearContents()
# Row 100
$ws.Range("H100").Value = 3565.6667
$ws.Range("I100").Value = 3512
$ws.Range("J100").Value = 3995
$ws.Range("K100").Value = 3512
$ws.Range("L100").Value = 3995
$ws.Range("M100").Value = -2971
$ws.Range("N100").Value = -5077
# Row 122
$ws.Range("H122").Value = 3836.9546
$ws.Range("I122").Value = 2612.6155
$ws.Range("J122").Value = 5605.4443
$ws.Range("K122").Value = 7837.8465
$ws.Range("L122").Value = 16816.3329
$ws.Range("M122").Value = -5387.8465
$ws.Range("N122").Value = -21716.3329
# Row 126
$ws.Range("H126").Value = 35716680
$ws.Range("I126").Value = 45456750
$ws.Range("J126").Value = 3093.3333
$ws.Range("K126").Value = 136370250
$ws.Range("L126").Value = 9279.999899999999
$ws.Range("M126").Value = -136367780
$ws.Range("N126").Value = -14219.9999
# Row 136
$ws.Range("H136").Value = 1963.5
$ws.Range("I136").Value = 1833.3226
$ws.Range("K136").Value = 5499.9678
$ws.Range("M136").Value = -2949.9678
$ws = $wb.Worksheets.Item("WVR")
# Row 4
$ws.Range("H4").Value = 1413674.5
$ws.Range("I4").Value = 2860300.2
$ws.Range("K4").Value = 2860300.2
$ws.Range("M4").Value = -2860187.2
# Row 43
$ws.Range("H43").Value = 25000
$ws.Range("I43").Value = 25000
$ws.Range("K43").Value = 25000
$ws.Range("M43").Value = -24851
# Row 122
$ws.Range("H122").Value = 1929
$ws.Range("I122").Value = 1799.6
$ws.Range("K122").Value = 5398.799999999999
$ws.Range("M122").Value = -2948.799999999999
